$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 227.83
$ws.Range("I15").Value = 227.83
$ws.Range("K15").Value = 683.49
$ws.Range("M15").Value = -514.49

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1750.25
$ws.Range("I43").Value = 2014.2858
$ws.Range("J43").Value = 1134.1666
$ws.Range("K43").Value = 2014.2858
$ws.Range("L43").Value = 1134.1666
$ws.Range("M43").Value = -1945.2858
$ws.Range("N43").Value = -1272.1666

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1059.16
$ws.Range("J112").Value = 1082.4584
$ws.Range("L112").Value = 3247.3752
$ws.Range("N112").Value = -5463.3752

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2647.7358
$ws.Range("J129").Value = 874.67566
$ws.Range("L129").Value = 2624.02698
$ws.Range("N129").Value = -12624.02698

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1264.5667
$ws.Range("I137").Value = 1126.0435
$ws.Range("J137").Value = 1719.7142
$ws.Range("K137").Value = 3378.1305
$ws.Range("L137").Value = 5159.142599999999
$ws.Range("M137").Value = -828.1305000000002
$ws.Range("N137").Value = -10259.1426

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24329.63
$ws.Range("I32").Value = 7301.1953
$ws.Range("K32").Value = 7301.1953
$ws.Range("M32").Value = -7014.1953

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 63469.06
$ws.Range("I45").Value = 79308.766
$ws.Range("K45").Value = 79308.766
$ws.Range("M45").Value = -78931.766

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2388.3333
$ws.Range("I61").Value = 1422.6154
$ws.Range("J61").Value = 3016.05
$ws.Range("K61").Value = 1422.6154
$ws.Range("L61").Value = 3016.05
$ws.Range("M61").Value = -1210.6154
$ws.Range("N61").Value = -3440.05

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1097
$ws.Range("I74").Value = 1031.0834
$ws.Range("J74").Value = 1294.75
$ws.Range("K74").Value = 1031.0834
$ws.Range("L74").Value = 1294.75
$ws.Range("M74").Value = -157.0834
$ws.Range("N74").Value = -3042.75

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1097
$ws.Range("I77").Value = 1031.0834
$ws.Range("J77").Value = 1294.75
$ws.Range("K77").Value = 5155.416999999999
$ws.Range("L77").Value = 6473.75
$ws.Range("M77").Value = -787.4169999999995
$ws.Range("N77").Value = -15209.75

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2126.2
$ws.Range("I122").Value = 1836.9333
$ws.Range("J122").Value = 2994
$ws.Range("K122").Value = 5510.7999
$ws.Range("L122").Value = 8982
$ws.Range("M122").Value = -3060.7999
$ws.Range("N122").Value = -13882

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 10857.17
$ws.Range("I132").Value = 11911.274
$ws.Range("K132").Value = 35733.822
$ws.Range("M132").Value = -33203.822

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2388.3333
$ws.Range("I136").Value = 1422.6154
$ws.Range("J136").Value = 3016.05
$ws.Range("K136").Value = 4267.8462
$ws.Range("L136").Value = 9048.150000000001
$ws.Range("M136").Value = -1717.8462
$ws.Range("N136").Value = -14148.15

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 333485700
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1880.9
$ws.Range("I58").Value = 1529.95
$ws.Range("K58").Value = 1529.95
$ws.Range("M58").Value = -1326.95

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2943.2307
$ws.Range("I99").Value = 2981.4
$ws.Range("J99").Value = 2919.375
$ws.Range("K99").Value = 2981.4
$ws.Range("L99").Value = 2919.375
$ws.Range("M99").Value = -1483.4
$ws.Range("N99").Value = -5915.375

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1108.2354
$ws.Range("I122").Value = 1096
$ws.Range("K122").Value = 3288
$ws.Range("M122").Value = -838

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2943.2307
$ws.Range("I126").Value = 2981.4
$ws.Range("J126").Value = 2919.375
$ws.Range("K126").Value = 8944.200000000001
$ws.Range("L126").Value = 8758.125
$ws.Range("M126").Value = -6474.200000000001
$ws.Range("N126").Value = -13698.125

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 83337096
$ws.Range("I132").Value = 90913570
$ws.Range("J132").Value = 71431210
$ws.Range("K132").Value = 272740710
$ws.Range("L132").Value = 214293630
$ws.Range("M132").Value = -272738180
$ws.Range("N132").Value = -214298690

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1331.0741
$ws.Range("I134").Value = 881
$ws.Range("J134").Value = 2400
$ws.Range("K134").Value = 2643
$ws.Range("L134").Value = 7200
$ws.Range("M134").Value = -108
$ws.Range("N134").Value = -12270

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1880.9
$ws.Range("I136").Value = 1529.95
$ws.Range("K136").Value = 4589.85
$ws.Range("M136").Value = -2039.85

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 51875.05
$ws.Range("I70").Value = 84558.414
$ws.Range("K70").Value = 253675.242
$ws.Range("M70").Value = -253360.242

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 51875.05
$ws.Range("I73").Value = 84558.414
$ws.Range("K73").Value = 253675.242
$ws.Range("M73").Value = -252583.242

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 776.2820400000001
$ws.Range("I107").Value = 681.625
$ws.Range("J107").Value = 800.70966
$ws.Range("K107").Value = 2044.875
$ws.Range("L107").Value = 2402.12898
$ws.Range("M107").Value = -124.875
$ws.Range("N107").Value = -6242.12898

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 815.35
$ws.Range("I131").Value = 372.8
$ws.Range("J131").Value = 864.5222
$ws.Range("K131").Value = 1118.4
$ws.Range("L131").Value = 2593.5666
$ws.Range("M131").Value = 3921.6
$ws.Range("N131").Value = -12673.5666

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 14982555
$ws.Range("I137").Value = 64321.188
$ws.Range("J137").Value = 33343458
$ws.Range("K137").Value = 192963.564
$ws.Range("L137").Value = 100030374
$ws.Range("M137").Value = -187863.564
$ws.Range("N137").Value = -100040574

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 12658.889
$ws.Range("I141").Value = 14848.571
$ws.Range("J141").Value = 4995
$ws.Range("K141").Value = 44545.713
$ws.Range("L141").Value = 14985
$ws.Range("M141").Value = -39365.713
$ws.Range("N141").Value = -25345

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2387.3
$ws.Range("I132").Value = 1752.6786
$ws.Range("J132").Value = 3868.0833
$ws.Range("K132").Value = 5258.0358
$ws.Range("L132").Value = 11604.2499
$ws.Range("M132").Value = -2728.0358
$ws.Range("N132").Value = -16664.2499

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2947.8333
$ws.Range("I136").Value = 2146.9285
$ws.Range("K136").Value = 6440.7855
$ws.Range("M136").Value = -3890.7855

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 586.2778
$ws.Range("I113").Value = 321.08334
$ws.Range("K113").Value = 963.2500200000001
$ws.Range("M113").Value = 1206.74998

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2478.6597
$ws.Range("I132").Value = 2830.6667
$ws.Range("J132").Value = 1857.4706
$ws.Range("K132").Value = 8492.000100000001
$ws.Range("L132").Value = 5572.4118
$ws.Range("M132").Value = -5962.000100000001
$ws.Range("N132").Value = -10632.4118

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2602.88
$ws.Range("I136").Value = 1006.2
$ws.Range("J136").Value = 3667.3333
$ws.Range("K136").Value = 3018.6
$ws.Range("L136").Value = 11001.9999
$ws.Range("M136").Value = -468.6000000000004
$ws.Range("N136").Value = -16101.9999
